# Auto-generated edit script: applies numeric corrections to Leve profit sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 118166.664
$ws.Range("J3").Value = 118166.664
$ws.Range("L3").Value = 118166.664
$ws.Range("N3").Value = -118394.664
$ws.Range("H15").Value = 1389964
$ws.Range("I15").Value = 1389964
$ws.Range("K15").Value = 4169892
$ws.Range("M15").Value = -4169723
$ws.Range("H86").Value = 4330.3076
$ws.Range("J86").Value = 4499.4443
$ws.Range("L86").Value = 4499.4443
$ws.Range("N86").Value = -6745.4443
$ws.Range("H89").Value = 4330.3076
$ws.Range("J89").Value = 4499.4443
$ws.Range("L89").Value = 22497.2215
$ws.Range("N89").Value = -33729.2215
$ws.Range("H98").Value = 2880.853
$ws.Range("I98").Value = 598.56665
$ws.Range("K98").Value = 598.56665
$ws.Range("M98").Value = 899.43335
$ws.Range("H99").Value = 2362.6667
$ws.Range("I99").Value = 1619.7693
$ws.Range("K99").Value = 4859.3079
$ws.Range("M99").Value = -3361.3079
$ws.Range("H102").Value = 118166.664
$ws.Range("J102").Value = 118166.664
$ws.Range("L102").Value = 118166.664
$ws.Range("N102").Value = -124656.664
$ws.Range("H121").Value = 1999.5
$ws.Range("J121").Value = 1999.5
$ws.Range("L121").Value = 5998.5
$ws.Range("N121").Value = -9492.5
$ws.Range("H122").Value = 2880.853
$ws.Range("I122").Value = 598.56665
$ws.Range("K122").Value = 1795.69995
$ws.Range("M122").Value = 654.3000500000001
$ws.Range("H138").Value = 2757.0908
$ws.Range("I138").Value = 2061.1
$ws.Range("J138").Value = 3827.8462
$ws.Range("K138").Value = 6183.299999999999
$ws.Range("L138").Value = 11483.5386
$ws.Range("M138").Value = -1043.299999999999
$ws.Range("N138").Value = -21763.5386
$ws.Range("H139").Value = 99959.75
$ws.Range("J139").Value = 99959.75
$ws.Range("L139").Value = 99959.75
$ws.Range("N139").Value = -110239.75

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 19150.525
$ws.Range("I32").Value = 19753.256
$ws.Range("K32").Value = 19753.256
$ws.Range("M32").Value = -19466.256
$ws.Range("H61").Value = 2553.6272
$ws.Range("I61").Value = 1818.2291
$ws.Range("K61").Value = 1818.2291
$ws.Range("M61").Value = -1606.2291
$ws.Range("H97").Value = 2529.7104
$ws.Range("I97").Value = 2053.4644
$ws.Range("J97").Value = 3863.2
$ws.Range("K97").Value = 2053.4644
$ws.Range("L97").Value = 3863.2
$ws.Range("M97").Value = -1557.4644
$ws.Range("N97").Value = -4855.2
$ws.Range("H110").Value = 1326.7778
$ws.Range("I110").Value = 1404
$ws.Range("K110").Value = 1404
$ws.Range("M110").Value = 641
$ws.Range("H122").Value = 2161.4
$ws.Range("I122").Value = 2106.4783
$ws.Range("J122").Value = 2341.8572
$ws.Range("K122").Value = 6319.4349
$ws.Range("L122").Value = 7025.571599999999
$ws.Range("M122").Value = -3869.4349
$ws.Range("N122").Value = -11925.5716
$ws.Range("H132").Value = 1234.8667
$ws.Range("I132").Value = 1183.862
$ws.Range("K132").Value = 3551.586
$ws.Range("M132").Value = -1021.586
$ws.Range("H136").Value = 2553.6272
$ws.Range("I136").Value = 1818.2291
$ws.Range("K136").Value = 5454.6873
$ws.Range("M136").Value = -2904.6873

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1295.7
$ws.Range("I94").Value = 1331.5952
$ws.Range("J94").Value = 1107.25
$ws.Range("K94").Value = 1331.5952
$ws.Range("L94").Value = 1107.25
$ws.Range("M94").Value = -880.5952
$ws.Range("N94").Value = -2009.25
$ws.Range("H107").Value = 16964.912
$ws.Range("I107").Value = 18844.068
$ws.Range("K107").Value = 18844.068
$ws.Range("M107").Value = -16924.068

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H43").Value = 21598
$ws.Range("J43").Value = 21598
$ws.Range("L43").Value = 21598
$ws.Range("N43").Value = -21966
$ws.Range("H88").Value = 33744.25
$ws.Range("J88").Value = 33744.25
$ws.Range("L88").Value = 33744.25
$ws.Range("N88").Value = -34556.25
$ws.Range("H91").Value = 33744.25
$ws.Range("J91").Value = 33744.25
$ws.Range("L91").Value = 33744.25
$ws.Range("N91").Value = -36552.25
$ws.Range("H101").Value = 21598
$ws.Range("J101").Value = 21598
$ws.Range("L101").Value = 21598
$ws.Range("N101").Value = -28088
$ws.Range("H132").Value = 13666.8545
$ws.Range("I132").Value = 13666.8545
$ws.Range("K132").Value = 41000.5635
$ws.Range("M132").Value = -38470.5635
$ws.Range("H134").Value = 2718
$ws.Range("I134").Value = 2051.32
$ws.Range("K134").Value = 6153.960000000001
$ws.Range("M134").Value = -3618.960000000001
$ws.Range("H140").Value = 112036.5
$ws.Range("J140").Value = 112036.5
$ws.Range("L140").Value = 112036.5
$ws.Range("N140").Value = -122396.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 10139.88
$ws.Range("I99").Value = 4749.5
$ws.Range("J99").Value = 11842.105
$ws.Range("K99").Value = 14248.5
$ws.Range("L99").Value = 35526.315
$ws.Range("M99").Value = -12002.5
$ws.Range("N99").Value = -40018.315

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 9999
$ws.Range("J23").Value = 9999
$ws.Range("L23").Value = 9999
$ws.Range("N23").Value = -10445
$ws.Range("H43").Value = 4297.5
$ws.Range("I43").Value = 2323.7273
$ws.Range("J43").Value = 26009
$ws.Range("K43").Value = 2323.7273
$ws.Range("L43").Value = 26009
$ws.Range("M43").Value = -2172.7273
$ws.Range("N43").Value = -26311
$ws.Range("H102").Value = 23972.348
$ws.Range("I102").Value = 28281.63
$ws.Range("K102").Value = 28281.63
$ws.Range("M102").Value = -26659.63
$ws.Range("H132").Value = 2720.5625
$ws.Range("I132").Value = 2735.2666
$ws.Range("K132").Value = 8205.799800000001
$ws.Range("M132").Value = -5675.799800000001
$ws.Range("H135").Value = 121233.5
$ws.Range("J135").Value = 131409.72
$ws.Range("L135").Value = 131409.72
$ws.Range("N135").Value = -141549.72
$ws.Range("H136").Value = 9460.655000000001
$ws.Range("J136").Value = 9460.655000000001
$ws.Range("L136").Value = 28381.965
$ws.Range("N136").Value = -33481.965
$ws.Range("H139").Value = 80795.78
$ws.Range("J139").Value = 80795.78
$ws.Range("L139").Value = 80795.78
$ws.Range("N139").Value = -91075.78
$ws.Range("H140").Value = 142978.5
$ws.Range("J140").Value = 170638
$ws.Range("L140").Value = 170638
$ws.Range("N140").Value = -180998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H29").Value = 31333.334
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 31333.334
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 31333.334
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -31923.334
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("H122").Value = 16817.273
$ws.Range("I122").Value = 16817.273
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 50451.819
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -48001.819
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 5398.815
$ws.Range("I132").Value = 3485
$ws.Range("J132").Value = 20709.334
$ws.Range("K132").Value = 10455
$ws.Range("L132").Value = 62128.00199999999
$ws.Range("M132").Value = -7925
$ws.Range("N132").Value = -67188.00199999999
$ws.Range("H136").Value = 3931.889
$ws.Range("I136").Value = 3799.7144
$ws.Range("J136").Value = 4394.5
$ws.Range("K136").Value = 11399.1432
$ws.Range("L136").Value = 13183.5
$ws.Range("M136").Value = -8849.143199999999
$ws.Range("N136").Value = -18283.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8444.923000000001
$ws.Range("I62").Value = 9224.75
$ws.Range("J62").Value = 8098.3335
$ws.Range("K62").Value = 9224.75
$ws.Range("L62").Value = 8098.3335
$ws.Range("M62").Value = -8600.75
$ws.Range("N62").Value = -9346.333500000001
$ws.Range("H65").Value = 8444.923000000001
$ws.Range("I65").Value = 9224.75
$ws.Range("J65").Value = 8098.3335
$ws.Range("K65").Value = 46123.75
$ws.Range("L65").Value = 40491.6675
$ws.Range("M65").Value = -43003.75
$ws.Range("N65").Value = -46731.6675
$ws.Range("H81").Value = 6319.6
$ws.Range("I81").Value = 6319.6
$ws.Range("K81").Value = 12639.2
$ws.Range("M81").Value = -11578.2
$ws.Range("H84").Value = 6319.6
$ws.Range("I84").Value = 6319.6
$ws.Range("K84").Value = 63196
$ws.Range("M84").Value = -57892
$ws.Range("H122").Value = 14325303
$ws.Range("I122").Value = 18568510
$ws.Range("K122").Value = 55705530
$ws.Range("M122").Value = -55703080
$ws.Range("H132").Value = 4561307
$ws.Range("I132").Value = 5017036
$ws.Range("J132").Value = 4017.4
$ws.Range("K132").Value = 15051108
$ws.Range("L132").Value = 12052.2
$ws.Range("M132").Value = -15048578
$ws.Range("N132").Value = -17112.2
$ws.Range("H139").Value = 124886.2
$ws.Range("J139").Value = 124886.2
$ws.Range("L139").Value = 124886.2
$ws.Range("N139").Value = -135166.2
